$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Collapse the three CORE COMPETENCIES detail paragraphs into one
#    short summary paragraph.
# ------------------------------------------------------------------
$bullet = [char]0x2022
$core = $d.Paragraphs.Item(6)
$core.Range.Text = "Statistical Analysis & Machine Learning " + $bullet + " Big Data & Data Engineering " + $bullet + " Data Visualization & Reporting"

# Remove the two now-redundant detail paragraphs that followed it.
$p7 = $d.Paragraphs.Item(7)
$p8 = $d.Paragraphs.Item(8)
$killRange = $d.Range($p7.Range.Start, $p8.Range.End)
$killRange.Delete()

# ------------------------------------------------------------------
# 2. Append a new "TECHNICAL SKILLS" section at the end of the
#    document, ahead of the final section properties.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$heading = $d.Paragraphs.Item($d.Paragraphs.Count)
$heading.Range.Text = "TECHNICAL SKILLS"
$heading.Style = "Heading 2"

$heading.Range.InsertParagraphAfter()
$line1 = $d.Paragraphs.Item($d.Paragraphs.Count)
$line1.Style = "Normal"
$line1.Range.Text = "STATISTICAL ANALYSIS & MACHINE LEARNING Advanced Statistical Modeling; Predictive Analytics; Data Mining; Machine Learning; Statistical Computing; A/B Testing; Meta-analytical Techniques"

$line1.Range.InsertParagraphAfter()
$line2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$line2.Style = "Normal"
$line2.Range.Text = "BIG DATA & DATA ENGINEERING Big Data Processing; Data Warehousing; Cloud Platforms; Databases; Data Governance; Streaming Data; Data Pipeline Optimization"

$line2.Range.InsertParagraphAfter()
$line3 = $d.Paragraphs.Item($d.Paragraphs.Count)
$line3.Style = "Normal"
$line3.Range.Text = "DATA VISUALIZATION & REPORTING Data Visualization; Geospatial Analysis; Interactive Dashboards; Statistical Reporting; Choropleths and Hexagonal Grid Maps for Demographic Visualization; Business Intelligence; Client Presentation"
